$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells G1 and H1
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the formatting of the existing header (F1) onto the new headers
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Fill elapsed time / cpu values for rows 2 through 10
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 7).Value = 1.127317944850074
    $ws.Cells.Item($row, 8).Value = 0.985
}
